# Updated MCH102 to MCH251
# Adds the MCH137-1 "CUBA ANTI-APARTHEID PUBLICATIONS" series row to Sheet1,
# restyles it to match the rest of the sheet, restores the frozen header
# pane/selection, and normalizes row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (row 2) ---------------------------------------------------
$cells = @("A2", "C2", "D2", "E2", "F2", "G2", "H2")
foreach ($addr in $cells) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Calibri"
    $c.Font.Size = 10
    $c.Font.ThemeColor = 1
}

$ws.Range("A2").Value = "MCH137-1"
$ws.Range("C2").Value = "CUBA ANTI-APARTHEID PUBLICATIONS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 21M | GRAP COUNT NUMER: NONE"

# --- Row heights -------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75

# --- Keep the header row frozen and select the new data row ------------------
$ws.Activate()
$ws.Range("A2:H2").Select()
$excel.ActiveWindow.FreezePanes = $true
